$wb = $excel.ActiveWorkbook
$wsColumns = $wb.Worksheets.Item("Columns")
$wsReports = $wb.Worksheets.Item("Reports")

# --- Table1 ("Columns" sheet): insert a new "Wrap Text" column before "Number Format" ---
$lo = $wsColumns.ListObjects.Item("Table1")

# Grow the table by one column (new column is always appended at the end by this engine)
$lo.Resize($wsColumns.Range("A1:K8"))

# Move the existing "Number Format" column (currently J) into the new last column (K)
$numberFormatValues = $wsColumns.Range("J1:J8").Value2
$wsColumns.Range("K1:K8").Value2 = $numberFormatValues
$wsColumns.Range("J1:J8").ClearContents()

# Populate the now-empty J column as the new "Wrap Text" column
$wsColumns.Range("J1").Value2 = "Wrap Text"
$wsColumns.Range("J4").Value2 = "Y"

# Carry the old "Number Format" column width (J) over to its new home (K), and give
# the new "Wrap Text" column (J) the same width as its neighbour ("Font Name", I)
$wsColumns.Columns.Item(11).ColumnWidth = $wsColumns.Columns.Item(10).ColumnWidth
$wsColumns.Columns.Item(10).ColumnWidth = $wsColumns.Columns.Item(9).ColumnWidth

# --- Column Width update for the "Title" row (row 4) ---
$wsColumns.Range("D4").Value2 = 12

# --- Active sheet / selection changes ---
# Previously "Reports" was the active tab with N3 selected; now "Columns" is active with D5
# selected, and "Reports" keeps a selection of A2 for when it is revisited.
$wsReports.Activate() | Out-Null
$wsReports.Range("A2").Select() | Out-Null
$wsColumns.Activate() | Out-Null
$wsColumns.Range("D5").Select() | Out-Null
